# Weekly Fruta/Hortaliza update: insert a new week's worth of "Piña / Caramelo"
# price rows (Especial, Primera, Segunda, Tercera) at the top of the existing
# "Caramelo" block (row 952), pushing the rest of the table down by 4 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 952 - everything from 952 downward shifts to 956+,
# and the sheet's used range grows from T1012 to T1016.
$ws.Rows("952:955").Insert()

# Common (constant-across-the-block) column values, copied from the
# neighbouring rows of this Piña/Mercado Mayorista Lo Valledor block.
$mercadoId = 6
$mercado = "Mercado Mayorista Lo Valledor de Santiago"
$region = "Metropolitana"
$codreg = 13
$tipo = "Fruta"
$productoId = 100108
$producto = "Tropicales y subtropicales"
$categoriaId = 100108005
$categoria = "Piña"
$variedad = "Caramelo"
$origen = "Ecuador"
$fecha = 44585

# Especial / Primera / Segunda / Tercera rows for the new reporting date.
$rows = @(
    @{ Row = 952; Calidad = "Especial"; Volumen = 108;  PMin = 13000; PMax = 14000; PProm = 13500; Unidad = "$/caja 10 unidades"; PrecioKg = 1350; KgUnidad = 10 },
    @{ Row = 953; Calidad = "Primera";  Volumen = 108;  PMin = 13000; PMax = 14000; PProm = 13500; Unidad = "$/caja 12 unidades"; PrecioKg = 1125; KgUnidad = 12 },
    @{ Row = 954; Calidad = "Segunda";  Volumen = 108;  PMin = 13000; PMax = 14000; PProm = 13500; Unidad = "$/caja 14 unidades"; PrecioKg = 964;  KgUnidad = 14 },
    @{ Row = 955; Calidad = "Tercera";  Volumen = 108;  PMin = 13000; PMax = 14000; PProm = 13500; Unidad = "$/caja 16 unidades"; PrecioKg = 844;  KgUnidad = 16 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
